# The deck had an accidental duplicate slide: slide 8 ("User Features:" /
# parking-capacity rating bullets) is a byte-for-byte duplicate of slide 9
# (which additionally carries the "Product Features (Continued)" title).
# Fix: remove the duplicate (slide 8); everything after it shifts up.
$p = $ppt.ActivePresentation
$p.Slides.Item(8).Delete()
